$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Program_sheet: rename the LMS program entry and drop the "PlayWright" row
# ---------------------------------------------------------------------------
$program = $wb.Worksheets.Item("Program_sheet")

# row 2: "LMSProgram" / "LMS" -> "LMSPrograms" / "LMSTeach"
$program.Range("A2").Value = "LMSPrograms"
$program.Range("B2").Value = "LMSTeach"

# row 8 ("PlayWright" / "Active") is no longer needed - remove its contents
$program.Range("A8:C8").Clear()

# ---------------------------------------------------------------------------
# Batch: rename the "InformationTechPath" program to "MobileSeleniumAuto"
# (and flag the renamed cells in red), add a log line for row 6, fix a
# couple of class counts, and rename "ITPath" -> "selenium"
# ---------------------------------------------------------------------------
$batch = $wb.Worksheets.Item("Batch")

$batchProgramCells = @("B2", "B3", "B4", "B5", "B7", "B8", "B9", "B10")
foreach ($cellRef in $batchProgramCells) {
    $cell = $batch.Range($cellRef)
    $cell.Value = "MobileSeleniumAuto"
    $cell.Font.Color = 255
}

# added log line: row 6 now also carries a (blank) Program Name cell,
# styled the same as the other renamed cells
$batch.Range("B6").Font.Color = 255

# fix Number of Classes values
$batch.Range("C2").Value = 9
$batch.Range("C4").Value = 11
$batch.Range("C5").Value = 12

# rename search helper value
$searchCell = $batch.Range("C12")
$searchCell.Value = "selenium"
$searchCell.Font.Color = 255
